$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume change (column E) values
# as scraped by the GitHub Actions job on Tue Sep 12 09:10:48 UTC 2023.

$ws.Cells.Item(2, 4).Value = "25.858.93"
$ws.Cells.Item(2, 5).Value = "  -0.34%  "
$ws.Cells.Item(3, 4).Value = "1.583.03"
$ws.Cells.Item(3, 5).Value = "  -2.03%  "
$ws.Cells.Item(4, 5).Value = "  -0.16%  "
$ws.Cells.Item(5, 4).Value = "209.94"
$ws.Cells.Item(5, 5).Value = "  -0.68%  "
$ws.Cells.Item(6, 5).Value = "  -0.17%  "
$ws.Cells.Item(7, 4).Value = "0.475"
$ws.Cells.Item(7, 5).Value = "  -2.63%  "
$ws.Cells.Item(8, 5).Value = "  -0.27%  "
$ws.Cells.Item(9, 4).Value = "0.0612"
$ws.Cells.Item(9, 5).Value = "  -1.27%  "
$ws.Cells.Item(10, 5).Value = "  -0.71%  "
$ws.Cells.Item(11, 5).Value = "  -0.29%  "
$ws.Cells.Item(12, 4).Value = "1.803.93"
$ws.Cells.Item(12, 5).Value = "  -1.98%  "
$ws.Cells.Item(13, 4).Value = "1.580.47"
$ws.Cells.Item(13, 5).Value = "  -2.16%  "
$ws.Cells.Item(14, 5).Value = "  -1.59%  "
$ws.Cells.Item(15, 4).Value = "0.503"
$ws.Cells.Item(15, 5).Value = "  -2.77%  "
$ws.Cells.Item(16, 4).Value = "25.861.13"
$ws.Cells.Item(16, 5).Value = "  -0.36%  "
$ws.Cells.Item(17, 4).Value = "0.0₃0722"
$ws.Cells.Item(17, 5).Value = "  -1.22%  "
$ws.Cells.Item(18, 5).Value = "  -2.54%  "
$ws.Cells.Item(19, 5).Value = "  -0.11%  "
$ws.Cells.Item(20, 4).Value = "192.95"
$ws.Cells.Item(20, 5).Value = "  +0.92%  "
$ws.Cells.Item(21, 4).Value = "4.18"
$ws.Cells.Item(21, 5).Value = "  -0.90%  "
$ws.Cells.Item(22, 4).Value = "9.33"
$ws.Cells.Item(22, 5).Value = "  -0.84%  "
$ws.Cells.Item(23, 4).Value = "5.92"
$ws.Cells.Item(23, 5).Value = "  -1.26%  "
$ws.Cells.Item(24, 5).Value = "  +0.23%  "
$ws.Cells.Item(25, 4).Value = "140.74"
$ws.Cells.Item(25, 5).Value = "  -1.59%  "
$ws.Cells.Item(26, 5).Value = "  -0.11%  "
$ws.Cells.Item(27, 4).Value = "1.69"
$ws.Cells.Item(27, 5).Value = "  -1.80%  "
$ws.Cells.Item(28, 4).Value = "15.05"
$ws.Cells.Item(28, 5).Value = "  -0.34%  "
$ws.Cells.Item(29, 4).Value = "6.42"
$ws.Cells.Item(29, 5).Value = "  -2.63%  "
$ws.Cells.Item(30, 5).Value = "  -4.63%  "
$ws.Cells.Item(31, 5).Value = "  -0.30%  "
$ws.Cells.Item(32, 5).Value = "  +0.50%  "
$ws.Cells.Item(33, 4).Value = "3.01"
$ws.Cells.Item(33, 5).Value = "  -1.98%  "
$ws.Cells.Item(34, 5).Value = "  +0.81%  "
$ws.Cells.Item(35, 5).Value = "  -2.09%  "
$ws.Cells.Item(36, 4).Value = "1.095.48"
$ws.Cells.Item(36, 5).Value = "  -2.40%  "
$ws.Cells.Item(37, 5).Value = "  -0.28%  "
$ws.Cells.Item(38, 5).Value = "  -1.80%  "
$ws.Cells.Item(39, 5).Value = "  -0.99%  "
$ws.Cells.Item(40, 4).Value = "0.500"
$ws.Cells.Item(40, 5).Value = "  -2.58%  "
$ws.Cells.Item(41, 4).Value = "0.775"
$ws.Cells.Item(41, 5).Value = "  -5.16%  "
$ws.Cells.Item(42, 4).Value = "0.802"
$ws.Cells.Item(42, 5).Value = "  +6.55%  "
$ws.Cells.Item(43, 4).Value = "93.06"
$ws.Cells.Item(43, 5).Value = "  -4.15%  "
$ws.Cells.Item(44, 4).Value = "5.10"
$ws.Cells.Item(44, 5).Value = "  +0.77%  "
$ws.Cells.Item(45, 4).Value = "1.717.90"
$ws.Cells.Item(45, 5).Value = "  -1.95%  "
$ws.Cells.Item(46, 4).Value = "0.0₆0111"
$ws.Cells.Item(46, 5).Value = "  -2.36%  "
$ws.Cells.Item(47, 4).Value = "1.50"
$ws.Cells.Item(47, 5).Value = "  +1.32%  "
$ws.Cells.Item(48, 4).Value = "53.03"
$ws.Cells.Item(48, 5).Value = "  -1.29%  "
$ws.Cells.Item(49, 5).Value = "  -1.53%  "
$ws.Cells.Item(50, 5).Value = "  -0.83%  "
$ws.Cells.Item(51, 5).Value = "  -0.18%  "

